$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells in the order that mirrors the target shared-string insertion
# order: reason, Account Blocked by FIA, reference, date, test
$ws.Range("C1").Value = "reason"
$ws.Range("C2").Value = "Account Blocked by FIA"
$ws.Range("D1").Value = "reference"
$ws.Range("E1").Value = "date"
$ws.Range("D2").Value = "test"

# Update remaining data cells
$ws.Range("A2").Value = 1008784483
$ws.Range("E2").Value = 20230907

# Set column C width to fit new content
$ws.Columns.Item(3).ColumnWidth = 21

# Update selection to reflect active cell E2
$ws.Range("E2").Select()
